$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 1).Value = 21.27851028916892
$ws.Cells.Item(2, 2).Value = 15.55683831460628
$ws.Cells.Item(2, 3).Value = 26.42125989957814
$ws.Cells.Item(3, 1).Value = 20.41054398523
$ws.Cells.Item(3, 2).Value = 12.55971118876687
$ws.Cells.Item(3, 3).Value = 28.99418136998837
$ws.Cells.Item(4, 1).Value = 20.48175051692525
$ws.Cells.Item(4, 2).Value = 14.18977926189117
$ws.Cells.Item(4, 3).Value = 27.64582283757333
$ws.Cells.Item(5, 1).Value = 24.90533879033199
$ws.Cells.Item(5, 2).Value = 15.57380600711907
$ws.Cells.Item(5, 3).Value = 33.72910077575387
$ws.Cells.Item(6, 1).Value = 26.14371771818085
$ws.Cells.Item(6, 2).Value = 18.784123780606
$ws.Cells.Item(6, 3).Value = 33.81043031600917
$ws.Cells.Item(7, 1).Value = 19.33553917934039
$ws.Cells.Item(7, 2).Value = 14.89028513076683
$ws.Cells.Item(7, 3).Value = 24.32551876465234
$ws.Cells.Item(8, 1).Value = 10.10494786249335
$ws.Cells.Item(8, 2).Value = 4.383035274927962
$ws.Cells.Item(8, 3).Value = 16.02588343053751
$ws.Cells.Item(9, 1).Value = 24.36126516973622
$ws.Cells.Item(9, 2).Value = 18.32602359876933
$ws.Cells.Item(9, 3).Value = 30.84520422215541
$ws.Cells.Item(10, 1).Value = 10.17521119104675
$ws.Cells.Item(10, 2).Value = 5.153487853980469
$ws.Cells.Item(10, 3).Value = 17.15549315441622
$ws.Cells.Item(11, 1).Value = 18.1259506470778
$ws.Cells.Item(11, 2).Value = 12.30201517077239
$ws.Cells.Item(11, 3).Value = 23.80353978360296
$ws.Cells.Item(12, 1).Value = 26.89618997383307
$ws.Cells.Item(12, 2).Value = 19.87579329386602
$ws.Cells.Item(12, 3).Value = 34.5491322018315
$ws.Cells.Item(13, 1).Value = 12.3731121221789
$ws.Cells.Item(13, 2).Value = 6.672386542026537
$ws.Cells.Item(13, 3).Value = 20.4214111117216
$ws.Cells.Item(14, 1).Value = 14.90595555650731
$ws.Cells.Item(14, 2).Value = 9.034286297458346
$ws.Cells.Item(14, 3).Value = 21.88222348580927
$ws.Cells.Item(15, 1).Value = 29.09134814925604
$ws.Cells.Item(15, 2).Value = 20.066136079619
$ws.Cells.Item(15, 3).Value = 39.81964163521254
$ws.Cells.Item(16, 1).Value = 26.15277604508424
$ws.Cells.Item(16, 2).Value = 19.00166449813064
$ws.Cells.Item(16, 3).Value = 33.72279404297961
$ws.Cells.Item(17, 1).Value = 27.01878608663255
$ws.Cells.Item(17, 2).Value = 20.00892014559372
$ws.Cells.Item(17, 3).Value = 34.37920474010433
$ws.Cells.Item(18, 1).Value = 20.451363408095
$ws.Cells.Item(18, 2).Value = 13.55979942888981
$ws.Cells.Item(18, 3).Value = 28.35570837133776
$ws.Cells.Item(19, 1).Value = 18.48415737433077
$ws.Cells.Item(19, 2).Value = 12.54777407637926
$ws.Cells.Item(19, 3).Value = 25.68848774323577
$ws.Cells.Item(20, 1).Value = 29.78967671897496
$ws.Cells.Item(20, 2).Value = 21.81758836940899
$ws.Cells.Item(20, 3).Value = 37.68716377843243
$ws.Cells.Item(21, 1).Value = 29.00352172621084
$ws.Cells.Item(21, 2).Value = 20.72137582728896
$ws.Cells.Item(21, 3).Value = 39.13441474356576
$ws.Cells.Item(22, 1).Value = 10.32008479892086
$ws.Cells.Item(22, 2).Value = 4.498811011462661
$ws.Cells.Item(22, 3).Value = 17.84377609444305
$ws.Cells.Item(23, 1).Value = 30.7184456876663
$ws.Cells.Item(23, 2).Value = 21.09564213582235
$ws.Cells.Item(23, 3).Value = 40.98512602876113
$ws.Cells.Item(24, 1).Value = 26.15277604508424
$ws.Cells.Item(24, 2).Value = 19.00166449813064
$ws.Cells.Item(24, 3).Value = 33.72279404297961
$ws.Cells.Item(25, 1).Value = 22.45743139137369
$ws.Cells.Item(25, 2).Value = 15.18403227171122
$ws.Cells.Item(25, 3).Value = 29.9728781276117
$ws.Cells.Item(26, 1).Value = 29.08172895045648
$ws.Cells.Item(26, 2).Value = 21.48128667914418
$ws.Cells.Item(26, 3).Value = 39.79328583225615
$ws.Cells.Item(27, 1).Value = 12.8528107335122
$ws.Cells.Item(27, 2).Value = 6.366208663950633
$ws.Cells.Item(27, 3).Value = 22.35732919842456
$ws.Cells.Item(28, 1).Value = 20.21135332901322
$ws.Cells.Item(28, 2).Value = 15.46136001639763
$ws.Cells.Item(28, 3).Value = 25.6784378468932
$ws.Cells.Item(29, 1).Value = 32.4907176783593
$ws.Cells.Item(29, 2).Value = 22.78856659037615
$ws.Cells.Item(29, 3).Value = 45.9409199178819
$ws.Cells.Item(30, 1).Value = 25.63160429653583
$ws.Cells.Item(30, 2).Value = 19.29800361316631
$ws.Cells.Item(30, 3).Value = 32.39483337421518
$ws.Cells.Item(31, 1).Value = 14.11299006136359
$ws.Cells.Item(31, 2).Value = 9.468160671614932
$ws.Cells.Item(31, 3).Value = 19.51187718578307
$ws.Cells.Item(32, 1).Value = 11.83841039793953
$ws.Cells.Item(32, 2).Value = 6.843350820844287
$ws.Cells.Item(32, 3).Value = 18.15829428916005
$ws.Cells.Item(33, 1).Value = 14.9126663498363
$ws.Cells.Item(33, 2).Value = 9.140409490426107
$ws.Cells.Item(33, 3).Value = 21.87543403599251
$ws.Cells.Item(34, 1).Value = 22.99966032146279
$ws.Cells.Item(34, 2).Value = 17.01734266543937
$ws.Cells.Item(34, 3).Value = 30.36593326389974
$ws.Cells.Item(35, 1).Value = 27.83191380176457
$ws.Cells.Item(35, 2).Value = 20.37563615980657
$ws.Cells.Item(35, 3).Value = 35.27962899662413
$ws.Cells.Item(36, 1).Value = 19.21681381218155
$ws.Cells.Item(36, 2).Value = 14.73727959639823
$ws.Cells.Item(36, 3).Value = 24.70863741987844
$ws.Cells.Item(37, 1).Value = 30.67830337994193
$ws.Cells.Item(37, 2).Value = 21.35295591150578
$ws.Cells.Item(37, 3).Value = 40.51957444987417
$ws.Cells.Item(38, 1).Value = 15.06629100716182
$ws.Cells.Item(38, 2).Value = 9.597162558172283
$ws.Cells.Item(38, 3).Value = 22.61564609847169
$ws.Cells.Item(39, 1).Value = 27.44616261027919
$ws.Cells.Item(39, 2).Value = 20.82179637369245
$ws.Cells.Item(39, 3).Value = 34.98637275057037
$ws.Cells.Item(40, 1).Value = 25.0341272051176
$ws.Cells.Item(40, 2).Value = 15.94980399383392
$ws.Cells.Item(40, 3).Value = 34.06441578255161
$ws.Cells.Item(41, 1).Value = 18.12467299610538
$ws.Cells.Item(41, 2).Value = 14.03729076394551
$ws.Cells.Item(41, 3).Value = 22.65718633788621
$ws.Cells.Item(42, 1).Value = 18.05405175132057
$ws.Cells.Item(42, 2).Value = 12.45314309420709
$ws.Cells.Item(42, 3).Value = 24.30188812837373
$ws.Cells.Item(43, 1).Value = 25.62376798132133
$ws.Cells.Item(43, 2).Value = 19.27880598353141
$ws.Cells.Item(43, 3).Value = 32.359567653879
$ws.Cells.Item(44, 1).Value = 15.83751587215647
$ws.Cells.Item(44, 2).Value = 12.37576138235381
$ws.Cells.Item(44, 3).Value = 19.776225454618
$ws.Cells.Item(45, 1).Value = 19.11362442153503
$ws.Cells.Item(45, 2).Value = 13.56084307645277
$ws.Cells.Item(45, 3).Value = 26.96463864426342
$ws.Cells.Item(46, 1).Value = 13.82464680342717
$ws.Cells.Item(46, 2).Value = 9.794669095079069
$ws.Cells.Item(46, 3).Value = 19.34444368245246
$ws.Cells.Item(47, 1).Value = 10.67978661526331
$ws.Cells.Item(47, 2).Value = 4.548455154323506
$ws.Cells.Item(47, 3).Value = 19.26959472406309
$ws.Cells.Item(48, 1).Value = 26.93256872099688
$ws.Cells.Item(48, 2).Value = 19.8519015347351
$ws.Cells.Item(48, 3).Value = 34.40756792008711
$ws.Cells.Item(49, 1).Value = 22.85890316152321
$ws.Cells.Item(49, 2).Value = 15.4958078097452
$ws.Cells.Item(49, 3).Value = 31.43089091686594
$ws.Cells.Item(50, 1).Value = 30.44638089347023
$ws.Cells.Item(50, 2).Value = 21.060365972208
$ws.Cells.Item(50, 3).Value = 40.37207822577744
$ws.Cells.Item(51, 1).Value = 20.27286043664231
$ws.Cells.Item(51, 2).Value = 12.86254066618309
$ws.Cells.Item(51, 3).Value = 28.75777648677797
$ws.Cells.Item(52, 1).Value = 13.56449833703462
$ws.Cells.Item(52, 2).Value = 9.409689999962877
$ws.Cells.Item(52, 3).Value = 18.20917938181911
$ws.Cells.Item(53, 1).Value = 27.60136401762282
$ws.Cells.Item(53, 2).Value = 21.02033799843995
$ws.Cells.Item(53, 3).Value = 34.79312091034125
$ws.Cells.Item(54, 1).Value = 29.14636431374536
$ws.Cells.Item(54, 2).Value = 20.03495732441929
$ws.Cells.Item(54, 3).Value = 40.35363914951202
$ws.Cells.Item(55, 1).Value = 29.2882427189699
$ws.Cells.Item(55, 2).Value = 21.25126233084593
$ws.Cells.Item(55, 3).Value = 39.56866907565718
$ws.Cells.Item(56, 1).Value = 13.98800543761862
$ws.Cells.Item(56, 2).Value = 8.947002003382984
$ws.Cells.Item(56, 3).Value = 20.02712308552932
$ws.Cells.Item(57, 1).Value = 11.65924032262891
$ws.Cells.Item(57, 2).Value = 6.433172718648613
$ws.Cells.Item(57, 3).Value = 18.14211700349625
$ws.Cells.Item(58, 1).Value = 30.64940940306926
$ws.Cells.Item(58, 2).Value = 21.70027910292186
$ws.Cells.Item(58, 3).Value = 40.58657261388718
$ws.Cells.Item(59, 1).Value = 17.14496354569453
$ws.Cells.Item(59, 2).Value = 12.47483167243192
$ws.Cells.Item(59, 3).Value = 22.62389549066331
$ws.Cells.Item(60, 1).Value = 30.78066905074913
$ws.Cells.Item(60, 2).Value = 20.78281140156295
$ws.Cells.Item(60, 3).Value = 41.22750849418951
$ws.Cells.Item(61, 1).Value = 14.88887767298378
$ws.Cells.Item(61, 2).Value = 8.956208354713194
$ws.Cells.Item(61, 3).Value = 21.5291395500169
$ws.Cells.Item(62, 1).Value = 18.15744150789238
$ws.Cells.Item(62, 2).Value = 12.2219941712868
$ws.Cells.Item(62, 3).Value = 23.75494338694534
$ws.Cells.Item(63, 1).Value = 15.51865398730565
$ws.Cells.Item(63, 2).Value = 9.163914493173362
$ws.Cells.Item(63, 3).Value = 24.26763918575239
$ws.Cells.Item(64, 1).Value = 30.80614208825699
$ws.Cells.Item(64, 2).Value = 20.98210853597818
$ws.Cells.Item(64, 3).Value = 41.47477835015054
$ws.Cells.Item(65, 1).Value = 11.64327464967391
$ws.Cells.Item(65, 2).Value = 6.517848784635598
$ws.Cells.Item(65, 3).Value = 17.43128550329988
$ws.Cells.Item(66, 1).Value = 22.83227520420498
$ws.Cells.Item(66, 2).Value = 16.54607797094488
$ws.Cells.Item(66, 3).Value = 30.25611285310313
$ws.Cells.Item(67, 1).Value = 21.79938783336812
$ws.Cells.Item(67, 2).Value = 16.27595422422494
$ws.Cells.Item(67, 3).Value = 27.89719036215173
$ws.Cells.Item(68, 1).Value = 27.02215455973724
$ws.Cells.Item(68, 2).Value = 19.24379904401171
$ws.Cells.Item(68, 3).Value = 35.66485360800263
$ws.Cells.Item(69, 1).Value = 10.93453536212381
$ws.Cells.Item(69, 2).Value = 5.323063163768153
$ws.Cells.Item(69, 3).Value = 19.79054120260843
$ws.Cells.Item(70, 1).Value = 24.85816931912767
$ws.Cells.Item(70, 2).Value = 16.30821597230861
$ws.Cells.Item(70, 3).Value = 34.5592510836908
$ws.Cells.Item(71, 1).Value = 19.76938926921078
$ws.Cells.Item(71, 2).Value = 14.550934979648
$ws.Cells.Item(71, 3).Value = 25.56935208302299
$ws.Cells.Item(72, 1).Value = 22.95537583904342
$ws.Cells.Item(72, 2).Value = 17.39938096330207
$ws.Cells.Item(72, 3).Value = 28.47600063475871
